$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.842.82'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '1.640.19'
$ws.Range("E3").Value = '  +0.58%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '''215.51'
$ws.Range("E5").Value = '  -0.06%  '

$ws.Range("D6").Value = '''0.5062'
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '''0.2590'
$ws.Range("E8").Value = '  +0.55%  '

$ws.Range("D9").Value = '''0.06439'
$ws.Range("E9").Value = '  +1.82%  '

$ws.Range("D10").Value = '''20.58'
$ws.Range("E10").Value = '  +5.70%  '

$ws.Range("D11").Value = '''0.07824'
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").Value = '''4.288'
$ws.Range("E12").Value = '  +0.97%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.866.64'
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.634.55'
$ws.Range("E14").Value = '  +0.14%  '

$ws.Range("D15").Value = '''0.5615'
$ws.Range("E15").Value = '  +2.33%  '

$ws.Range("E16").Value = '  +0.88%  '

$ws.Range("D17").Value = '''63.25'
$ws.Range("E17").Value = '  -0.64%  '

$ws.Range("D18").Value = '25.866.69'
$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").Value = '''1.001'
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("D20").Value = '''193.52'
$ws.Range("E20").Value = '  -0.31%  '

$ws.Range("D21").Value = '''4.380'
$ws.Range("E21").Value = '  -0.79%  '

$ws.Range("D22").Value = '''9.953'
$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("D23").Value = '''6.141'
$ws.Range("E23").Value = '  +1.79%  '

$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '''1.793'
$ws.Range("E25").Value = '  -6.35%  '

$ws.Range("D26").Value = '''140.38'
$ws.Range("E26").Value = '  -1.16%  '

$ws.Range("D27").Value = '''0.1239'
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").Value = '''6.835'
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("D29").Value = '''15.58'
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = '''1.244'
$ws.Range("E30").Value = '  +0.53%  '

$ws.Range("D31").Value = '''0.04998'
$ws.Range("E31").Value = '  +2.61%  '

$ws.Range("D32").Value = '''3.302'
$ws.Range("E32").Value = '  +2.04%  '

$ws.Range("D33").Value = '''3.241'
$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("D34").Value = '''1.581'
$ws.Range("E34").Value = '  +2.49%  '

$ws.Range("D35").Value = '''2.379'
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("D36").Value = '''0.9071'
$ws.Range("E36").Value = '  +1.35%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.5595'
$ws.Range("E37").Value = '  +1.37%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '''2.571'
$ws.Range("E38").Value = '  +1.27%  '

$ws.Range("D39").Value = '1.126.32'
$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").Value = '''0.01570'
$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").Value = '''0.9958'
$ws.Range("E41").Value = '  -0.52%  '

$ws.Range("D42").Value = '''5.530'
$ws.Range("E42").Value = '  -0.77%  '

$ws.Range("D43").Value = '''0.8031'
$ws.Range("E43").Value = '  +0.84%  '

$ws.Range("D44").Value = '''98.53'
$ws.Range("E44").Value = '  +1.47%  '

$ws.Range("D45").Value = '1.778.38'
$ws.Range("E45").Value = '  +0.68%  '

$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").Value = '  -8.05%  '

$ws.Range("D47").Value = '''55.72'
$ws.Range("E47").Value = '  +1.75%  '

$ws.Range("D48").Value = '''0.4268'
$ws.Range("E48").Value = '  -4.02%  '

$ws.Range("D49").Value = '''7.789'
$ws.Range("E49").Value = '  +3.40%  '

$ws.Range("D50").Value = '''0.05047'
$ws.Range("E50").Value = '  -1.70%  '

$ws.Range("D51").Value = '''0.9973'
$ws.Range("E51").Value = '  -0.61%  '

